$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.186135172843933
$ws.Range("B1").Value = 1.999900579452515
$ws.Range("C1").Value = 6.24809455871582
$ws.Range("D1").Value = 2.30059814453125
$ws.Range("E1").Value = 1.196532487869263
